$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8693733215332031
$ws.Range("B1").Value = 1.649717092514038
$ws.Range("C1").Value = 6.123432636260986
$ws.Range("D1").Value = 1.874566912651062
$ws.Range("E1").Value = 1.132223606109619
